$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = '67.718.61'
    "E2" = '  -1.23%  '
    "D3" = '3.781.77'
    "E3" = '  +0.50%  '
    "E4" = '  -0.05%  '
    "D5" = '''595.83'
    "E5" = '  +0.40%  '
    "D6" = '''167.34'
    "E6" = '  +0.02%  '
    "D7" = '3.779.28'
    "E7" = '  +0.48%  '
    "E8" = '  +0.05%  '
    "D9" = '''0.520'
    "E9" = '  -0.36%  '
    "E10" = '  -0.15%  '
    "D11" = '''6.28'
    "E11" = '  -2.02%  '
    "D12" = '''0.449'
    "E12" = '  -0.03%  '
    "D13" = '''0.0000254'
    "E13" = '  -2.56%  '
    "D14" = '''36.09'
    "E14" = '  -0.27%  '
    "D15" = '4.413.94'
    "E15" = '  +0.51%  '
    "D16" = '3.780.75'
    "E16" = '  +0.51%  '
    "D17" = '67.691.16'
    "E17" = '  -1.19%  '
    "D18" = '''18.40'
    "E18" = '  +2.30%  '
    "D19" = '''7.03'
    "E19" = '  +0.32%  '
    "E20" = '  -1.07%  '
    "E21" = '  -6.45%  '
    "D22" = '''458.28'
    "E22" = '  -1.52%  '
    "D23" = '''0.696'
    "E23" = '  -0.28%  '
    "D24" = '''0.0000154'
    "E24" = '  +3.87%  '
    "D25" = '''83.28'
    "E25" = '  -0.93%  '
    "D26" = '''11.97'
    "E26" = '  +0.03%  '
    "D27" = '''2.14'
    "E27" = '  -2.34%  '
    "D28" = '''10.04'
    "E28" = '  -0.21%  '
    "E29" = '  +0.14%  '
    "E30" = '  -0.31%  '
    "D31" = '''2.26'
    "E31" = '  +3.49%  '
    "D32" = '''7.23'
    "E32" = '  -1.64%  '
    "D33" = '''29.69'
    "E33" = '  -1.11%  '
    "D34" = '''9.12'
    "E34" = '  -0.74%  '
    "D35" = '''0.999'
    "E35" = '  +0.30%  '
    "B36" = 'Hedera'
    "C36" = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    "D36" = '''0.100'
    "E36" = '  -0.39%  '
    "B37" = 'dogwifhat'
    "C37" = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
    "D37" = '''3.38'
    "E37" = '  -0.85%  '
    "B38" = 'Kaspa'
    "C38" = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
    "D38" = '''0.138'
    "E38" = '  +0.43%  '
    "B39" = 'Mantle'
    "C39" = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    "D39" = '''0.995'
    "E39" = '  -0.69%  '
    "B40" = 'Filecoin'
    "C40" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "D40" = '''5.76'
    "E40" = '  -0.66%  '
    "B41" = 'FirstDigitalUSD'
    "C41" = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    "D41" = '''0.999'
    "E41" = '  -0.04%  '
    "B42" = 'USDe'
    "C42" = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
    "D42" = '''1.00'
    "E42" = '  +0.00%  '
    "B43" = 'Arweave'
    "C43" = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
    "D43" = '''45.83'
    "E43" = '  +4.20%  '
    "B44" = 'OKB'
    "C44" = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
    "D44" = '''48.14'
    "E44" = '  +2.90%  '
    "B45" = 'TheGraph'
    "C45" = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
    "D45" = '''0.299'
    "E45" = '  -1.25%  '
    "B46" = 'Monero'
    "C46" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "D46" = '''149.42'
    "E46" = '  +2.60%  '
    "B47" = 'Cosmos'
    "C47" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "D47" = '''8.33'
    "E47" = '  -1.99%  '
    "B48" = 'Bittensor'
    "C48" = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
    "D48" = '''393.83'
    "E48" = '  +1.12%  '
    "B49" = 'Stacks'
    "C49" = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
    "D49" = '''1.83'
    "E49" = '  -4.64%  '
    "B50" = 'EnergySwap'
    "C50" = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    "D50" = '''26.33'
    "E50" = '  +1.91%  '
    "B51" = 'Maker'
    "C51" = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    "D51" = '2.725.30'
    "E51" = '  -1.34%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
